# Fruta / hortaliza, semanal
# Insert this week's new price record as a new row right above the
# existing row 26, shifting all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new blank row at row 26 (pushes old rows 26..105 down to 27..106)
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new week's data
$ws.Cells.Item(26, 1).Value = 8
$ws.Cells.Item(26, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44414
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(26, 5).Value = 4
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = "Poroto verde"
$ws.Cells.Item(26, 8).Value = "Magnum"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 540
$ws.Cells.Item(26, 11).Value = 31000
$ws.Cells.Item(26, 12).Value = 32000
$ws.Cells.Item(26, 13).Value = 31500
$ws.Cells.Item(26, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Perú"
$ws.Cells.Item(26, 16).Value = 1260
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
